$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, D (or $null if unchanged), E
$updates = @(
    ,@(2, "28.512.18", "  -1.88%  ")
    ,@(3, "1.959.37", "  -0.28%  ")
    ,@(4, "1.011", "  +0.54%  ")
    ,@(5, "322.58", "  -1.49%  ")
    ,@(6, "1.011", "  +0.50%  ")
    ,@(7, "0.4802", "  -3.94%  ")
    ,@(8, "0.4074", "  -3.31%  ")
    ,@(9, "53.95", "  +1.70%  ")
    ,@(10, "0.08512", "  -7.43%  ")
    ,@(11, "1.061", "  -3.49%  ")
    ,@(12, "22.47", "  -1.83%  ")
    ,@(13, "1.965.12", "  -2.94%  ")
    ,@(14, "7.588", "  -3.54%  ")
    ,@(15, "6.172", "  -4.16%  ")
    ,@(16, "1.014", "  +0.62%  ")
    ,@(17, "90.75", "  -0.70%  ")
    ,@(18, "0.00001074", "  -2.46%  ")
    ,@(19, "0.06620", "  -0.92%  ")
    ,@(20, "18.47", "  -3.91%  ")
    ,@(21, "1.011", "  +0.66%  ")
    ,@(22, "5.855", "  -1.68%  ")
    ,@(23, "28.530.10", "  -1.93%  ")
    ,@(24, $null, "  -5.38%  ")
    ,@(25, $null, "  +0.51%  ")
    ,@(26, "2.195.50", "  -0.85%  ")
    ,@(27, "156.47", "  +0.00%  ")
    ,@(28, "20.34", "  -1.29%  ")
    ,@(29, "2.178", "  -3.65%  ")
    ,@(30, "5.846", "  -5.61%  ")
    ,@(31, $null, "  -1.79%  ")
    ,@(32, "0.9880", "  -5.23%  ")
    ,@(33, "0.09666", "  -1.90%  ")
    ,@(34, "1.458", "  -4.54%  ")
    ,@(35, "5.647", "  -2.26%  ")
    ,@(36, "3.695", "  +0.40%  ")
    ,@(37, "9.120", "  +1.87%  ")
    ,@(38, $null, "  -3.57%  ")
    ,@(39, "0.06199", "  -1.89%  ")
    ,@(40, "1.255", "  -3.50%  ")
    ,@(41, "0.6244", "  -3.09%  ")
    ,@(42, "11.21", "  -1.99%  ")
    ,@(43, "1.011", "  +0.65%  ")
    ,@(44, "0.1919", "  -3.59%  ")
    ,@(45, "1.357", "  +5.42%  ")
    ,@(46, "0.5963", "  -4.45%  ")
    ,@(47, "13.05", "  -1.87%  ")
    ,@(48, "2.064", "  -5.34%  ")
    ,@(49, "3.412", "  -1.56%  ")
    ,@(50, "0.06818", "  -1.58%  ")
    ,@(51, "111.25", "  -1.35%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]
    if ($null -ne $dVal) {
        $dCell = $ws.Range("D$row")
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
        $dCell.ClearFormats()
    }
    $ws.Range("E$row").Value = $eVal
}
